$d = $word.ActiveDocument

# --- Step 1: add the date run "24.01.2026" to the last (empty) Body Text paragraph ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$dateRange = $lastPara.Range
$dateRange.InsertAfter("24.01.2026")
$dateRange.LanguageID = "en-US"

# --- Step 2: insert a new empty Body Text paragraph, then a Body Text paragraph
#             with the "Today I learned..." text (with spell-check markers around
#             each "css" occurrence), both appended after the paragraph above. ---
$insertionPoint = $d.Range($d.Content.End, $d.Content.End)

$newParagraphsXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
  '<w:pPr><w:pStyle w:val="BodyText"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '</w:p>' +
  '<w:p>' +
  '<w:pPr><w:pStyle w:val="BodyText"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Today I learned about </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>css</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> selectors, </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>css</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> grid, flex-box and also margins, padding and a lot more. I designed the hero section with the given hero-image also I learned how to create utility classes and how they are very useful for designing a website fast. Also how to use alignment for text. Basically today I created three major section of our growth webpage. Which are hero-section, testimonial section and pricing section. </w:t></w:r>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData>' +
  '</pkg:part>' +
  '</pkg:package>'

$insertionPoint.InsertXML($newParagraphsXml)

Write-Output ("Paragraphs after edit: " + $d.Paragraphs.Count)
